$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The order in which brand-new shared strings are first introduced matters
# for reproducing the exact sharedStrings.xml layout, so cells are written
# in that specific sequence (row 13, then 8, then 7, then 10, then 14, with
# rows 9/11 - which introduce no new text - filled in alongside).

# Row 13 (new): Order / Defensive / Debt to Equity / debttoequity / asc
$ws.Range("A13").Value = "Order"
$ws.Range("B13").Value = "Defensive"
$ws.Range("C13").Value = "Debt to Equity"
$ws.Range("D13").Value = "debttoequity"
$ws.Range("E13").Value = "asc"

# Row 8 (repurposed): View / Defensive / Price to Book / pricetobook
$ws.Range("A8").Value = "View"
$ws.Range("B8").Value = "Defensive"
$ws.Range("C8").Value = "Price to Book"
$ws.Range("D8").Value = "pricetobook"
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

# Row 9 (new): View / Risky / Price to Book / pricetobook
$ws.Range("A9").Value = "View"
$ws.Range("B9").Value = "Risky"
$ws.Range("C9").Value = "Price to Book"
$ws.Range("D9").Value = "pricetobook"

# Row 7 (new): Search / Risky / Price to Earnings / pricetoearnings / ~gt~ / 50
$ws.Range("A7").Value = "Search"
$ws.Range("B7").Value = "Risky"
$ws.Range("C7").Value = "Price to Earnings"
$ws.Range("D7").Value = "pricetoearnings"
$ws.Range("E7").Value = "~gt~"
$ws.Range("F7").Value = 50

# Row 10 (new): View / Risky / Price to Revenue / pricetorevenue
$ws.Range("A10").Value = "View"
$ws.Range("B10").Value = "Risky"
$ws.Range("C10").Value = "Price to Revenue"
$ws.Range("D10").Value = "pricetorevenue"

# Row 11 (old row 8 data moved down): Search / Standard / Market Cap / marketcap / ~gt~ / 2000000000
$ws.Range("A11").Value = "Search"
$ws.Range("B11").Value = "Standard"
$ws.Range("C11").Value = "Market Cap"
$ws.Range("D11").Value = "marketcap"
$ws.Range("E11").Value = "~gt~"
$ws.Range("F11").Value = 2000000000

# Row 14 (new): Order / Risky / Revenue Growth / revenuegrowth / desc
$ws.Range("A14").Value = "Order"
$ws.Range("B14").Value = "Risky"
$ws.Range("C14").Value = "Revenue Growth"
$ws.Range("D14").Value = "revenuegrowth"
$ws.Range("E14").Value = "desc"

$ws.Range("F13").Select()
